# Celerio-generated "address.xlsx" update:
# add a new "search_full_text" search-criteria row to the "Search" sheet,
# inserted above the existing "streetName"/"city" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search")

# Insert a new row above the current row 4 (the streetName row); this
# shifts the streetName row to 5 and the city row to 6.
$ws.Rows.Item(4).Insert()

# Fill in the new row with the search_full_text label/placeholder pair.
$ws.Range("A4").Value = "`${msg.getProperty('search_full_text')}"
$ws.Range("B4").Value = "`${search_full_text}"

Write-Output "A4=$($ws.Range('A4').Value())"
Write-Output "B4=$($ws.Range('B4').Value())"
Write-Output "A5=$($ws.Range('A5').Value())"
Write-Output "B5=$($ws.Range('B5').Value())"
Write-Output "A6=$($ws.Range('A6').Value())"
Write-Output "B6=$($ws.Range('B6').Value())"
